$d = $word.ActiveDocument

# --- helper: identity "find/replace" that forces Word to rewrite the run(s)
# spanning the given exact text, stripping any <w:proofErr/> markers that
# sat inside the span and coalescing the runs it touches into one run.
function Reflow-Text($text) {
    $d.Content.Find.Execute($text, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $text, 2) | Out-Null
}

# 1) "Zoals b.v.: Google font Lato: https://www.google.com/fonts"
#    -> proofErr marks removed; text itself unchanged.
Reflow-Text "Zoals b.v.: Google font Lato: https://www.google.com/fonts"

# 2) "Wishlist: Eigen lijst ..."
Reflow-Text "Wishlist: Eigen lijst bij kunnen houden van films die je graag wilt hebben"

# 3) "De web app moet goed werken ..."
Reflow-Text "De web app moet goed werken in zo veel mogelijk webbrowsers en moet er tevens goed uitzien en functioneren op mobiele apparaten zoals een smartphone en tablet"

# 4) "Projecteisen (communicatie, aanleveren content, budget, oplevering etc.)"
Reflow-Text "Projecteisen (communicatie, aanleveren content, budget, oplevering etc.)"

# 5) "Aanleveren content: ..." -- real text rewrite + bookmark relocation
#    First drop the existing _GoBack bookmark (Word keeps a single one, it
#    moves to track the last edit position).
$existingGoBack = $d.Bookmarks("_GoBack")
$existingGoBack.Delete()

$d.Content.Find.Execute( `
    "Aanleveren content: De content van de web app kan getest en na de projectperiode ondergebracht worden bij het webhosting pakket van Menno van der Krift", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Aanleveren content: De web app kan getest worden op, en later ondergebracht worden onder het webhosting pakket van Menno van der Krift", `
    2) | Out-Null

# Re-insert _GoBack right after "Aanleveren content: De web "
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.StartsWith("Aanleveren content: De web app kan getest")) {
        $mark = $para.Range.Start + ("Aanleveren content: De web ").Length
        $bmRange = $d.Range($mark, $mark)
        $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
        break
    }
}

# 6) "Budget: €0,-"
Reflow-Text "Budget: €0,-"

# 7) "Oplevering: ... Tevens kan de web app bij het webhost pakket van Menno van der Krift worden ondergebracht"
Reflow-Text "Tevens kan de web app bij het webhost pakket van Menno van der Krift worden ondergebracht"
